$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill in the daily power records that were missing for 2018-12-12 and
# --- 2018-12-13 (rows 124/125 already existed as blank calculated-column
# --- rows in the table; only Date/Start/End were empty).
$ws.Range("A124").Value = 43446
$ws.Range("B124").Value = 0
$ws.Range("C124").Value = 0

$ws.Range("A125").Value = 43447
$ws.Range("B125").Value = 0
$ws.Range("C125").Value = 0

# --- Add a brand new row 126 to the sheet. Inserting (rather than just
# --- writing into a blank row) makes the new row inherit the number
# --- formatting of the row above it, same as Excel does when a table
# --- grows by one row. The inserted row initially spans all 6 columns;
# --- clear the A:C cells back out since the new record has no Date/Start/
# --- End value yet (matches the table's calculated-columns-only row).
$ws.Rows(126).Insert()
$ws.Range("A126:C126").Clear()

# --- Grow the table so the new row becomes part of it.
$tbl = $ws.ListObjects.Item(1)
$tbl.Resize($ws.Range("A1:F126"))

# --- Re-enter the calculated-column formulas across the existing data rows;
# --- writing them as one range fill groups them into shared formulas,
# --- matching how Excel stores an AutoFilled calculated table column.
$ws.Range("D117:D125").Formula = "=(C117-B117)* 1440"
$ws.Range("E117:E125").Formula = "=IF(C117>B117, (C117-B117)*1440, (B117-C117)*1440)"
$ws.Range("F117:F125").Formula = "=ABS((C117-B117)*1440)"

# --- New row 126 gets its own (non-shared) calculated-column formulas.
$ws.Range("D126").Formula = "=(C126-B126)* 1440"
$ws.Range("E126").Formula = "=IF(C126>B126, (C126-B126)*1440, (B126-C126)*1440)"
$ws.Range("F126").Formula = "=ABS((C126-B126)*1440)"

# --- Move the selection down to the newly added row, as Excel would leave
# --- it after typing the new record.
$ws.Range("A126").Select()
